# Applies the "example instead of custom domain" edit:
#   1. Refreshes the cached "Date Placeholder" text (auto date field) from
#      4/28/21 -> 5/2/21 on every slide layout, the slide master, and the
#      notes master.
#   2. Updates two copy-text runs on slide 1:
#        "Automate A/B, A/B/n, Canary, and Conformance experiments"
#          -> "Automate A/B(/n), Canary, and Conformance experiments"
#        "Safely promote winning version "
#          -> "Find and promote winning version "

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

$oldDate = "4/28/21"
$newDate = "5/2/21"

# Slide master
$master = $p.SlideMaster
Update-DateShapes $master.Shapes $oldDate $newDate

# Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes $layout.Shapes $oldDate $newDate
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes $oldDate $newDate

# Slide 1 copy updates
$slide = $p.Slides.Item(1)

$rectAB = $slide.Shapes.Item("Rounded Rectangle 53")
$rectAB.TextFrame.TextRange.Text = "Automate A/B(/n), Canary, and Conformance experiments"

$rectPromote = $slide.Shapes.Item("Rounded Rectangle 54")
$tr = $rectPromote.TextFrame.TextRange
$firstRunLen = ("Safely promote winning version ").Length
$tr.Characters(1, $firstRunLen).Text = "Find and promote winning version "
